$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.674.37"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.928.36"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "327.10"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").Value = "0.4828"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "0.4062"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "0.08194"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "1.010"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").Value = "23.77"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.070"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.889.31"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "7.292"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "91.44"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "0.06873"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "17.65"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "29.643.99"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "5.647"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "11.99"
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("D24").Value = "2.203"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "2.129.97"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").Value = "156.47"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "6.375"
$ws.Range("E27").Value = "  -4.81%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "2.091"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "120.97"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "1.007"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("D32").Value = "0.09607"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "5.614"
$ws.Range("D34").Value = "3.557"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "1.394"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "0.06530"
$ws.Range("E36").Value = "  +6.12%  "
$ws.Range("D37").Value = "0.02281"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").Value = "1.210"
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("D39").Value = "0.5931"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "10.77"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").Value = "7.869"
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("D42").Value = "2.532"
$ws.Range("E42").Value = "  +2.98%  "
$ws.Range("D43").Value = "0.1845"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "1.280"
$ws.Range("D45").Value = "0.07526"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").Value = "12.32"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "0.5554"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "1.960"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "117.77"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "2.424"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "71.90"
$ws.Range("E51").Value = "  -1.45%  "
